$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, matching the style used by existing
# header cells in row 1 (bold/centered/bordered style, s="1" -> Range.Style "Normal" won't
# reproduce custom xf, so copy format from an existing header cell instead).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style (font/border/alignment) of an existing header cell (H1) onto the
# newly added header cells so they match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows: column I ("I0") and column J ("IF") values.
$iValues = @{
  2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
  11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 7; 18 = 7
}
$jValues = @{
  2 = 4; 3 = 5; 4 = 4; 5 = 6; 6 = 6; 7 = 5; 8 = 5; 9 = 4; 10 = 5;
  11 = 5; 12 = 5; 13 = 7; 14 = 5; 15 = 5; 16 = 4; 17 = 9; 18 = 8
}

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
